$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 4 (old row4 -> row5, old row5 -> row6) ---
$ws.Rows("4").Insert()
$ws.Range("A4:I4").ClearFormats()

# --- 2. Row 1: new title text, merge A1:I1, centered/wrapped bold title ---
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Telavi Municipality"
$ws.Range("A1:I1").Merge()
$ws.Rows("1").RowHeight = 51
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true

# --- 3. Row 2: "(End of year, persons)" label - row height only change ---
$ws.Rows("2").RowHeight = 14.5

# --- 4. Row 3: year header row unaffected, only A3 font cosmetics ---
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# --- 5. Row 4 (new): "family with disabilities Persons " with values ---
$ws.Range("A4").Value = "family with disabilities Persons "
$cols = @("B","C","D","E","F","G","H","I")
$values4 = @(953,897,888,898,896,900,882,878)
for ($i = 0; $i -lt 8; $i++) {
    $ws.Range($cols[$i] + "4").Value = $values4[$i]
}
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true
$ws.Range("A4").Borders.Item(8).LineStyle = 1
$ws.Range("A4").Borders.Item(8).Weight = 2
$ws.Range("B4:I4").NumberFormat = "#\ ##0"
$ws.Range("B4:I4").Font.Name = "Arial"
$ws.Range("B4:I4").Font.Size = 10
$ws.Rows("4").RowHeight = 24.75

# --- 6. Row 5 (was row 4): relabel "disabilities Persons ", new values, drop top border ---
$ws.Range("A5").Value = "disabilities Persons "
$values5 = @(1053,984,967,977,980,985,957,958)
for ($i = 0; $i -lt 8; $i++) {
    $ws.Range($cols[$i] + "5").Value = $values5[$i]
}
$ws.Range("A5").Borders.Item(8).LineStyle = -4142
$ws.Range("B5:I5").NumberFormat = "#\ ##0"
$ws.Rows("5").RowHeight = 21

# --- 7. Row 6 (was row 5): Source row - drop top border on label cell only ---
$ws.Range("A6").Borders.Item(8).LineStyle = -4142
$ws.Rows("6").RowHeight = 27.75

# --- 8. Column width for A ---
$ws.Columns("A").ColumnWidth = 19.98

# --- 9. Selection on the new title range (matches authored selection state) ---
$ws.Range("A1:I1").Select()

Write-Output "done"
